# Update RF (raising factor) values in column I for the 2-RAP gear rows
# (2025 data / RF update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 gets its own updated RF value
$ws.Range("I16").Value = 26.84678571428572

# Rows 17 through 44 share the same updated RF value
$ws.Range("I17:I44").Value = 52.69357142857143
